$d = $word.ActiveDocument

# --- Edit 1: merge the "DB created and filled..." sentence back into a
# single run, removing the gramStart/gramEnd proofErr split around
# "diagrams". A same-text Find & Replace (format-preserving, match case)
# collapses the three runs into one clean run.
$oldFulfilled = "DB created and filled with dummy data, started sequence diagrams and uploaded files to GitLab."
$d.Content.Find.Execute($oldFulfilled, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $oldFulfilled, 2) | Out-Null

# --- Edit 2: fill in the previously-empty "Purpose of meeting" and
# "Fulfilled Tasks" cells for the 3/11 row.
$table = $d.Tables.Item(1)
$table.Cell(6, 2).Range.Text = "Finished Use Case and Class diagram, also completed two of three sequence diagrams. Have scheduled meeting for 16/11 at 11am (Online). Created issues on GitLab repository to plan work to be done before next meeting"
$table.Cell(6, 3).Range.Text = "All members to check and complete issues in GitLab repository and think of suitable third sequence diagram to complete for deadline 15/11 at 2pm."
